$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New expense row (19) goes right below the last existing entry (18) and
# above the blank buffer rows that precede the TOTAL row (22). Copy the
# formatting (date / currency / centered "y") from row 18 so the new row
# matches the rest of the table, then fill in the actual values.
$ws.Range("A18:E18").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A19").Value = 42494
$ws.Range("B19").Value = "Chain"
$ws.Range("C19").Value = "Amazon"
$ws.Range("D19").Value = 59.84
$ws.Range("E19").Value = "y"

# TOTAL (D22) is =SUM(D2:D21), so it recalculates automatically to include
# the new row.

$ws.Range("E20").Select()
